$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: section title
$ws.Range("D25").Value = "TestDePerformanceConNodoSocketsMismaLan"

# Row 26: header row (copy of row 6 / row 26 used to be blank style-only cells)
$ws.Range("D26").Value = 41107
$ws.Range("E26").Value = "Duración en milis"
$ws.Range("F26").Value = "Input(tareas/ms)"
$ws.Range("G26").Value = "Output(tareas/ms)"
$ws.Range("H26").Value = "O/I"
$ws.Range("I26").Value = "Msg/s(I)"
$ws.Range("J26").Value = "Msg/s(O)"

# Row 27
$ws.Range("D27").Value = "[1T->(wolfy)->1R]"
$ws.Range("E27").Value = 10000
$ws.Range("F27").Value = 15.975199999999999
$ws.Range("G27").Value = 0.98019999999999996
$ws.Range("H27").Formula = "=G27/F27"
$ws.Range("I27").Formula = "=F27*1000"
$ws.Range("J27").Formula = "=G27*1000"

# Row 28
$ws.Range("D28").Value = "[2T->(wolfy)->1R]"
$ws.Range("E28").Value = 10000
$ws.Range("F28").Value = 16.666799999999999
$ws.Range("G28").Value = 1.0943000000000001
$ws.Range("H28").Formula = "=G28/F28"
$ws.Range("I28").Formula = "=F28*1000"
$ws.Range("J28").Formula = "=G28*1000"

# Row 29
$ws.Range("D29").Value = "[4T->(wolfy)->1R]"
$ws.Range("E29").Value = 10000
$ws.Range("F29").Value = 17.142299999999999
$ws.Range("G29").Value = 1.0909
$ws.Range("H29").Formula = "=G29/F29"
$ws.Range("I29").Formula = "=F29*1000"
$ws.Range("J29").Formula = "=G29*1000"

# Row 30
$ws.Range("D30").Value = "[8T->(wolfy)->1R]"
$ws.Range("E30").Value = 10016
$ws.Range("F30").Value = 16.459265175718802
$ws.Range("G30").Value = 1.14007587859424
$ws.Range("H30").Formula = "=G30/F30"
$ws.Range("I30").Formula = "=F30*1000"
$ws.Range("J30").Formula = "=G30*1000"

# Row 31
$ws.Range("D31").Value = "[16T->(wolfy)->1R]"
$ws.Range("E31").Value = 10000
$ws.Range("F31").Value = 16.159199999999998
$ws.Range("G31").Value = 1.1767000000000001
$ws.Range("H31").Formula = "=G31/F31"
$ws.Range("I31").Formula = "=F31*1000"
$ws.Range("J31").Formula = "=G31*1000"

# Row 32
$ws.Range("D32").Value = "[32T->(wolfy)->1R]"
$ws.Range("E32").Value = 10000
$ws.Range("F32").Value = 15.631399999999999
$ws.Range("G32").Value = 1.1367
$ws.Range("H32").Formula = "=G32/F32"
$ws.Range("I32").Formula = "=F32*1000"
$ws.Range("J32").Formula = "=G32*1000"

# Row 33
$ws.Range("D33").Value = "[200T->(wolfy)->1R]"
$ws.Range("E33").Value = 10062
$ws.Range("F33").Value = 4.1171735241502603
$ws.Range("G33").Value = 1.27181474855893
$ws.Range("H33").Formula = "=G33/F33"
$ws.Range("I33").Formula = "=F33*1000"
$ws.Range("J33").Formula = "=G33*1000"

# Re-prioritize the conditional-formatting colour scales covering the rows
# that now contain data, so they bubble to the top priority (Excel does this
# automatically when a rule's range gets edited/recalculated).
$ws.Range("H27:H29,H33").FormatConditions.Item(1).SetFirstPriority()
$ws.Range("H32").FormatConditions.Item(1).SetFirstPriority()
$ws.Range("H31").FormatConditions.Item(1).SetFirstPriority()
$ws.Range("H30").FormatConditions.Item(1).SetFirstPriority()

$ws.Range("G17").Select()
